# "ortografia y solucios exceso de creditos"
#
# 1) "durante el periodo escolar 20-2" -> "... 18-1"              (ortografia)
# 2) Table "Creditos" column: row "1. Pin Pong"          1 -> 2   (exceso de creditos)
# 3) Table "Creditos" column: row "3. Taller de Lectura" 1 -> 2   (exceso de creditos)
# 4) Remove the rows "4. Baile Moderno" and "5. Regueton" entirely
#    (those activities' credits were not actually owed)

$d = $word.ActiveDocument

# --- 1) Fix the school period referenced in the body paragraph -------------
# Locate the paragraph that mentions "periodo escolar" and surgically replace
# just the "20-2" token it contains (keeps the surrounding runs untouched).
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $ptext = $p.Range.Text
    $idx = $ptext.IndexOf("periodo escolar 20-2")
    if ($idx -ge 0) {
        $tokenStart = $p.Range.Start + $idx + ("periodo escolar ").Length
        $tokenEnd = $tokenStart + ("20-2").Length
        $tokenRange = $d.Range($tokenStart, $tokenEnd)
        $tokenRange.Text = "18-1"
        break
    }
}

# --- 2)/3) Correct the credit counts that were under-reported --------------
$table = $d.Tables.Item(1)

function Set-CreditValueForRow($labelSubstring, $newValue) {
    for ($r = 1; $r -le $table.Rows.Count; $r++) {
        $labelCell = $table.Rows.Item($r).Cells.Item(1)
        if ($labelCell.Range.Text.Contains($labelSubstring)) {
            $creditCell = $table.Rows.Item($r).Cells.Item(3)
            $digitRange = $d.Range($creditCell.Range.Start, $creditCell.Range.Start + 1)
            $digitRange.Text = $newValue
            return
        }
    }
}

Set-CreditValueForRow "Pin Pong" "2"
Set-CreditValueForRow "Taller de Lectura" "2"

# --- 4) Delete the extra activity rows that caused the credit overage ------
# Find rows by their label text and delete them (highest index first so the
# remaining rows' indices stay valid while iterating).
function Remove-RowWithLabel($labelSubstring) {
    for ($r = $table.Rows.Count; $r -ge 1; $r--) {
        $labelCell = $table.Rows.Item($r).Cells.Item(1)
        if ($labelCell.Range.Text.Contains($labelSubstring)) {
            $table.Rows.Item($r).Delete()
            return
        }
    }
}

Remove-RowWithLabel "Regueton"
Remove-RowWithLabel "Baile Moderno"
